$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text while we write the
# new values, so Excel does not auto-convert numeric-looking strings (like
# "47.72") into real numbers. This preserves the original inline-string text
# cell semantics used throughout the sheet.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.660.92"
$ws.Range("E2").Value = "  +7.01%  "
$ws.Range("D3").Value = "1.944.97"
$ws.Range("E3").Value = "  +5.35%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "341.20"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "0.4780"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("D8").Value = "0.4117"
$ws.Range("E8").Value = "  +6.83%  "
$ws.Range("D9").Value = "47.72"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").Value = "1.031"
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("D12").Value = "22.64"
$ws.Range("E12").Value = "  +6.68%  "
$ws.Range("D13").Value = "1.935.93"
$ws.Range("E13").Value = "  +3.93%  "
$ws.Range("D14").Value = "6.136"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").Value = "7.352"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "91.67"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D18").Value = "0.00001056"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "0.06677"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "18.00"
$ws.Range("E20").Value = "  +4.26%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "29.628.97"
$ws.Range("E22").Value = "  +6.86%  "
$ws.Range("D23").Value = "5.583"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  +3.78%  "
$ws.Range("D25").Value = "2.288"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "2.171.04"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("D27").Value = "161.13"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "20.16"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "2.167"
$ws.Range("E29").Value = "  +5.10%  "
$ws.Range("D30").Value = "5.632"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("D31").Value = "122.87"
$ws.Range("E31").Value = "  +3.55%  "
$ws.Range("D32").Value = "1.007"
$ws.Range("E32").Value = "  +7.26%  "
$ws.Range("D33").Value = "0.09654"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "1.468"
$ws.Range("E34").Value = "  +10.93%  "
$ws.Range("D35").Value = "3.678"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("D36").Value = "5.482"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("D37").Value = "0.06256"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").Value = "0.02310"
$ws.Range("E38").Value = "  +4.61%  "
$ws.Range("D39").Value = "8.465"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").Value = "0.6064"
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("D42").Value = "10.68"
$ws.Range("E42").Value = "  +6.39%  "
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").Value = "0.1895"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "2.369"
$ws.Range("E46").Value = "  +32.21%  "
$ws.Range("D47").Value = "0.5706"
$ws.Range("E47").Value = "  +4.89%  "
$ws.Range("D48").Value = "12.46"
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("D49").Value = "0.07412"
$ws.Range("E49").Value = "  +8.42%  "
$ws.Range("D50").Value = "1.986"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").Value = "112.72"
$ws.Range("E51").Value = "  +1.86%  "

# Restore the normal (default) cell style so the written cells do not carry
# an explicit text-format style attribute, matching the original workbook.
$priceVolumeRange.Style = "Normal"
